$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row to match the latest snapshot.
# For Price cells whose new value looks like a plain number (e.g. "328.59"), force the cell
# to keep storing the value as text (as it was originally, t="inlineStr"/shared string) by
# temporarily applying a text number format, then resetting the style back to Normal so no
# stray formatting is left behind.

$ws.Range('D2').Value = '28.405.89'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '1.804.51'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '328.59'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9993'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4461'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.85%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3778'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +7.65%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '44.51'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.47%  '
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07495'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.58'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.001'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.626'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.28%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.292'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '1.801.75'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06798'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '80.51'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.51'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.316'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '28.395.41'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.415'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '20.45'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '154.04'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.346'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.75%  '
$ws.Range('D29').Value = '2.005.72'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '132.07'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.254'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.003'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.813'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.09362'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.2279'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.81%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.12'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06367'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02334'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6588'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.29%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.149'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  -3.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.118'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9991'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.78'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6081'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.809'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '128.25'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.031'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07091'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.153'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.67%  '
